$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Game" sheet: fill in Name / ESRB_ID / Online for the games that were
#    previously only stub rows (GameID 13..28), in the same top-to-bottom
#    order they appear in the sheet so the shared-string table comes out in
#    the same sequence as the real edit.
# ---------------------------------------------------------------------------
$game = $wb.Worksheets.Item("Game")

$game.Range("B14").Value = "Missile Command"
$game.Range("C14").Value = 1
$game.Range("E14").Value = 0

$game.Range("B15").Value = "Crash Bandicoot"
$game.Range("C15").Value = 5
$game.Range("E15").Value = 0

$game.Range("B16").Value = "Star Wars Battlefront"
$game.Range("C16").Value = 3
$game.Range("E16").Value = 0

$game.Range("B18").Value = "Halo Reach"
$game.Range("C18").Value = 4
$game.Range("E18").Value = 1

$game.Range("B19").Value = "Guitar Hero: On Tour"
$game.Range("C19").Value = 2
$game.Range("E19").Value = 0

$game.Range("B20").Value = "Monster Hunter World"
$game.Range("C20").Value = 3
$game.Range("E20").Value = 1

$game.Range("B21").Value = "Mario Kart Wii"
$game.Range("C21").Value = 1
$game.Range("E21").Value = 1

$game.Range("B22").Value = "Left 4 Dead"
$game.Range("C22").Value = 4
$game.Range("E22").Value = 1

$game.Range("B23").Value = "Skyrim "
$game.Range("C23").Value = 4
$game.Range("E23").Value = 0

$game.Range("B24").Value = "Destiny"
$game.Range("C24").Value = 3
$game.Range("E24").Value = 1

$game.Range("B25").Value = "Super Smash Bros. Brawl"
$game.Range("C25").Value = 3
$game.Range("E25").Value = 1

$game.Range("B26").Value = "Sonic The Hedgehog"
$game.Range("C26").Value = 1
$game.Range("E26").Value = 0

$game.Range("B27").Value = "Zelda Link's Awakening"
$game.Range("C27").Value = 1
$game.Range("E27").Value = 0

$game.Range("B28").Value = "Donkey Kong 64"
$game.Range("C28").Value = 1
$game.Range("E28").Value = 0

$game.Range("B29").Value = "Last of Us "
$game.Range("C29").Value = 4
$game.Range("E29").Value = 0

# Column B ("Name") got manually widened (and lost its auto best-fit flag)
# now that longer titles live in it.
$game.Columns.Item(2).ColumnWidth = 19.6

# ---------------------------------------------------------------------------
# 2. "ESRB" sheet: add the ESRB_Name lookup column (B), plus the ESRB_ID key
#    values (A) it lines up with, for the 5 rating tiers.
# ---------------------------------------------------------------------------
$esrb = $wb.Worksheets.Item("ESRB")

$esrb.Range("B1").Value = "ESRB_Name"

$esrb.Range("A2").Value = 1
$esrb.Range("B2").Value = "E"

$esrb.Range("A3").Value = 2
$esrb.Range("B3").Value = "E 10+"

$esrb.Range("A4").Value = 3
$esrb.Range("B4").Value = "T"

$esrb.Range("A5").Value = 4
$esrb.Range("B5").Value = "M"

$esrb.Range("A6").Value = 5
$esrb.Range("B6").Value = "KA"

$esrb.Columns.Item(2).ColumnWidth = 10.3

# ---------------------------------------------------------------------------
# 3. Move the active sheet/selection from "Console" to "Game" (set the
#    non-active sheets' lingering selection first, then activate+select on
#    "Game" last so it ends up as the tab that's actually shown).
# ---------------------------------------------------------------------------
$esrb.Activate()
$esrb.Range("C5").Select()

$game.Activate()
$game.Range("D17").Select()
